$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.377.08"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.19%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.937.37"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.12%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.40%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7713"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +6.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "246.13"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.91%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.35%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3202"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.77"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07039"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.37%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7808"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08030"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.937.00"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.08%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.64"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.43"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.64%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.368.33"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "255.75"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007932"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.73%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.787"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.192.01"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.35%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.36%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.726"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.43%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.552"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.72"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1347"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.07"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.268"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.371"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.517"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.417"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.43%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.112"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05160"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.95%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.281"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7463"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.47%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.785"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.67%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01952"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.811"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "78.44"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.415"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4503"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.971"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.003"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8364"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.01"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.754"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.502"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "985.61"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +10.79%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.14"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4145"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.25%  "
